$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 17, shifting rows 17:23 down to 18:24
$ws.Rows("17:17").Insert()

# Populate the new row 17 with its values
$ws.Range("A17").Value = 5
$ws.Range("B17").Value = "Macroferia Regional de Talca"
$ws.Range("C17").Value = "Maule"
$ws.Range("D17").Value = 44729
$ws.Range("D17").NumberFormat = $ws.Range("D18").NumberFormat
$ws.Range("E17").Value = 7
$ws.Range("F17").Value = 100112040
$ws.Range("G17").Value = "Cilantro"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 150
$ws.Range("K17").Value = 7000
$ws.Range("L17").Value = 7000
$ws.Range("M17").Value = 7000
$ws.Range("N17").Value = "$/caja 36 atados"
$ws.Range("O17").Value = "Región del Maule"
$ws.Range("P17").Value = 194
$ws.Range("Q17").Value = 36
$ws.Range("R17").Value = "Hortaliza"
